# Week 11 Software Testing.pptx - update stale "last edited" date stamps.
#
#   Plain-text date footers ("2561.11.06" -> "2562.10.25") live on:
#     - every slide's Date Placeholder (slides 2..63)
#     - every slide layout's Date Placeholder (all 11 layouts)
#     - the slide master's Date Placeholder
#
#   Field-backed date footers ("05/11/61" -> "25/10/62", a cached
#   datetimeFigureOut field) live on the Handout Master and Notes Master.
#
# The shape holding the date text is not always at the same index (it is
# usually #4, but drifts to #3 on a few slides), so we scan every shape's
# text instead of hard-coding an index/name.

$OLD_PLAIN  = "2561.11.06"
$NEW_PLAIN  = "2562.10.25"
$OLD_FIELD  = "05/11/61"
$NEW_FIELD  = "25/10/62"

function Update-ShapeDate {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }

        $tr = $shp.TextFrame.TextRange
        $txt = $tr.Text

        if ($txt -eq $OLD_PLAIN) {
            $tr.Text = $NEW_PLAIN
        } elseif ($txt -eq $OLD_FIELD) {
            $tr.Text = $NEW_FIELD
        }
    }
}

$p = $ppt.ActivePresentation

# 1) Every regular slide.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    Update-ShapeDate $slide.Shapes
}

# 2) The slide master.
$master = $p.SlideMaster
Update-ShapeDate $master.Shapes

# 3) Every slide layout hanging off the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-ShapeDate $layout.Shapes
}

# 4) Handout master / notes master date fields.
if ($p.HandoutMaster) {
    Update-ShapeDate $p.HandoutMaster.Shapes
}
if ($p.NotesMaster) {
    Update-ShapeDate $p.NotesMaster.Shapes
}
